$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.287.07'
$ws.Range('D3').Value = '2.511.73'
$ws.Range('E3').Value = '  -0.23%  '
$origStyle = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = $origStyle
$ws.Range('E4').Value = '  -0.07%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '540.07'
$ws.Range('D5').Style = $origStyle
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.53'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  -0.69%  '
$origStyle = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = $origStyle
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('D9').Value = '2.524.64'
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('E10').Value = '  +1.16%  '
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('E12').Value = '  -1.37%  '
$ws.Range('E13').Value = '  -1.94%  '
$ws.Range('D14').Value = '2.959.91'
$ws.Range('E14').Value = '  -0.18%  '
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.28'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('D16').Value = '59.182.82'
$ws.Range('E17').Value = '  -0.03%  '
$ws.Range('D18').Value = '2.526.44'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('E19').Value = '  +0.86%  '
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.30'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  +0.84%  '
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '326.47'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  +1.16%  '
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('E23').Value = '  +1.98%  '
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.44'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  +5.12%  '
$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.425'
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('E26').Value = '  +0.52%  '
$ws.Range('E27').Value = '  -0.16%  '
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.69'
$ws.Range('D28').Style = $origStyle
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.76'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0779'
$ws.Range('E30').Value = '  +1.05%  '
$ws.Range('E31').Value = '  +0.44%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '169.02'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  +3.53%  '
$ws.Range('E33').Value = '  +6.88%  '
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('E35').Value = '  +2.39%  '
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.55'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  +0.48%  '
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.15'
$ws.Range('D37').Style = $origStyle
$ws.Range('E37').Value = '  -1.77%  '
$ws.Range('E38').Value = '  -0.30%  '
$ws.Range('E39').Value = '  -0.37%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.826'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  +2.67%  '
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.66'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  +0.51%  '
$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '285.53'
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  +2.12%  '
$ws.Range('E43').Value = '  +1.44%  '
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.998'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  -0.06%  '
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '131.59'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  +7.41%  '
$ws.Range('E46').Value = '  +1.75%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.88'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  +0.08%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0934'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  +0.15%  '
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('E50').Value = '  -0.21%  '
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.53'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  -0.84%  '
